# HTTPS.xlsx GSC export refresh: the rolling window of daily "Pages" counts
# advances by one day (2025-11-19 .. 2026-02-16  ->  2025-11-20 .. 2026-02-17).
# Net effect on the "Chart" sheet's data table (rows 2-91, columns A:C):
#   - row 2 (the oldest day, 2025-11-19) drops off
#   - every remaining row's date label / Non-HTTPS / Pages data moves up one row
#   - a new row 91 is appended for the newest day, 2026-02-17 (Pages = 30)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the oldest day's row; Excel shifts rows 3:91 up into 2:90 for us,
# carrying each row's date label and Pages count along with it.
$ws.Rows.Item(2).Delete()

# Append the new day as row 91. Build the date text via a throwaway formula
# cell and paste only its resulting *value* back into A91 - this avoids
# Excel's literal-entry autodetection turning "2026-02-17" into a date
# serial number, so it lands as plain text like every other date cell in
# the column (and without leaving any extra/unused cell style behind).
$scratch = $ws.Cells.Item(200, 10)
$scratch.Formula = "=""2026-02-17"""
$scratch.Copy()
$ws.Cells.Item(91, 1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$scratch.Clear()
$excel.CutCopyMode = $false

$ws.Cells.Item(91, 2).Value = 0
$ws.Cells.Item(91, 3).Value = 30
